$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = "ơi, ạ, à, Vâng, nhé, À, Thôi, nhỉ, hả, Chao ôi"
$ws.Range("C12").Value = "VN, Nguyễn, Văn, Mỹ, Hùng, Hải, Việt, Nam, Hà Nội, Hoàng"
$ws.Range("C14").Value = "TP., TP, UBND, SV, ĐH, TB, DN, HIV, LĐ, HS"
$ws.Range("C17").Value = "cả, chính, ngay, thôi, rồi, cái, Ngay, thật, đâu, mà"
$ws.Range("C20").Value = "như thế, như vậy, làm sao, nhất là, thế nào, ngày càng, có lẽ, vì sao, một mình, Như vậy"
$ws.Range("C21").Value = "phó, viên, siêu, tổng, tái, bất, liên, hoá, vô, Phó"
